# "cambios y busqeda general"
# The TIPO column (L) on the "Table1" sheet held specific vehicle model
# names (EXPEDITION, CARGO, SPLINTER, PROMASTER, FRONTIER, HILUX, RAM 2500,
# GLC300) for each of the 8 data rows. They are all replaced with a single
# generic value, "jungla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 12).Value = "jungla"
}

# Update the current selection/navigation state left in the sheet: the
# active cell moves to F2 (single cell selection, no frozen/scrolled
# top-left cell override).
$ws.Range("F2").Select()
